{"js": "// The document contains a single table of simple addition/subtraction\n// \"fact practice\" problems (20 rows x 5 columns = 100 cells). This edit\n// regenerates the whole practice sheet: every cell's expression text is\n// replaced by a new one, in row-major (top-to-bottom, left-to-right)\n// order. The table's shape (20 rows, 5 columns) and all other document\n// content (the date line, fonts, paragraph/cell formatting) stay the\n// same \u2014 only the `w:t` text inside each cell's run changes.\nconst newValues = [\n  [\"64+13=\", \"28+12=\", \"86-46=\", \"40+11=\", \"95-42=\"],\n  [\"99-95=\", \"4+48=\", \"14+78=\", \"3+91=\", \"86-8=\"],\n  [\"8-0=\", \"6+42=\", \"13+65=\", \"32-12=\", \"12-10=\"],\n  [\"57+19=\", \"35+51=\", \"66-19=\", \"73-35=\", \"37+35=\"],\n  [\"30-24=\", \"20+29=\", \"94-16=\", \"39+10=\", \"26+71=\"],\n  [\"41-25=\", \"74-74=\", \"23+25=\", \"39-27=\", \"73-53=\"],\n  [\"87-81=\", \"11+77=\", \"55-33=\", \"6-1=\", \"39+3=\"],\n  [\"96-87=\", \"67-11=\", \"72+25=\", \"17+71=\", \"19+66=\"],\n  [\"46+33=\", \"32-19=\", \"43+46=\", \"90-31=\", \"21-0=\"],\n  [\"95-16=\", \"1+28=\", \"0+16=\", \"95-39=\", \"54-39=\"],\n  [\"39+36=\", \"41-37=\", \"76-57=\", \"99-81=\", \"6+22=\"],\n  [\"1+36=\", \"96-89=\", \"79-71=\", \"26+15=\", \"7+82=\"],\n  [\"95-52=\", \"28+22=\", \"3+15=\", \"57-32=\", \"67-21=\"],\n  [\"18+76=\", \"76-41=\", \"82-78=\", \"86-9=\", \"24+5=\"],\n  [\"99-13=\", \"17+66=\", \"12+54=\", \"68-37=\", \"24-22=\"],\n  [\"24-18=\", \"7+47=\", \"79+4=\", \"60+20=\", \"41+32=\"],\n  [\"81-42=\", \"92-47=\", \"3+83=\", \"33+43=\", \"64-51=\"],\n  [\"31-6=\", \"91-58=\", \"81-19=\", \"47-7=\", \"7+31=\"],\n  [\"85-11=\", \"37+60=\", \"89-4=\", \"91-28=\", \"66-9=\"],\n  [\"79-52=\", \"34-28=\", \"21+49=\", \"59-54=\", \"25+61=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of simple addition/subtraction\n# \"fact practice\" problems (20 rows x 5 columns = 100 cells). This edit\n# regenerates the whole practice sheet: every cell's expression text is\n# replaced by a new one, in row-major (top-to-bottom, left-to-right)\n# order. The table's shape (20 rows, 5 columns) and all other document\n# content (the date line, fonts, paragraph/cell formatting) stay the\n# same -- only the text inside each cell changes.\n\n$newValues = @(\n  @(\"64+13=\", \"28+12=\", \"86-46=\", \"40+11=\", \"95-42=\"),\n  @(\"99-95=\", \"4+48=\", \"14+78=\", \"3+91=\", \"86-8=\"),\n  @(\"8-0=\", \"6+42=\", \"13+65=\", \"32-12=\", \"12-10=\"),\n  @(\"57+19=\", \"35+51=\", \"66-19=\", \"73-35=\", \"37+35=\"),\n  @(\"30-24=\", \"20+29=\", \"94-16=\", \"39+10=\", \"26+71=\"),\n  @(\"41-25=\", \"74-74=\", \"23+25=\", \"39-27=\", \"73-53=\"),\n  @(\"87-81=\", \"11+77=\", \"55-33=\", \"6-1=\", \"39+3=\"),\n  @(\"96-87=\", \"67-11=\", \"72+25=\", \"17+71=\", \"19+66=\"),\n  @(\"46+33=\", \"32-19=\", \"43+46=\", \"90-31=\", \"21-0=\"),\n  @(\"95-16=\", \"1+28=\", \"0+16=\", \"95-39=\", \"54-39=\"),\n  @(\"39+36=\", \"41-37=\", \"76-57=\", \"99-81=\", \"6+22=\"),\n  @(\"1+36=\", \"96-89=\", \"79-71=\", \"26+15=\", \"7+82=\"),\n  @(\"95-52=\", \"28+22=\", \"3+15=\", \"57-32=\", \"67-21=\"),\n  @(\"18+76=\", \"76-41=\", \"82-78=\", \"86-9=\", \"24+5=\"),\n  @(\"99-13=\", \"17+66=\", \"12+54=\", \"68-37=\", \"24-22=\"),\n  @(\"24-18=\", \"7+47=\", \"79+4=\", \"60+20=\", \"41+32=\"),\n  @(\"81-42=\", \"92-47=\", \"3+83=\", \"33+43=\", \"64-51=\"),\n  @(\"31-6=\", \"91-58=\", \"81-19=\", \"47-7=\", \"7+31=\"),\n  @(\"85-11=\", \"37+60=\", \"89-4=\", \"91-28=\", \"66-9=\"),\n  @(\"79-52=\", \"34-28=\", \"21+49=\", \"59-54=\", \"25+61=\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n  for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n    $tbl.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n  }\n}\n"}
